$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update conversion note text in A1 ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$hoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 12.2 = 49268.29 pesos`n✅ 49268.29 pesos = 12.19 = 974.37 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update rate values ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 82
$tasas.Range("O10").Value = 4040
$tasas.Range("N12").Value = 4042.62
$tasas.Range("O12").Value = 79.95
